# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to match the refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetExhibition.Range("F2").Value = 4350
$sheetExhibition.Range("F11").Value = 157
$sheetExhibition.Range("F12").Value = 1622
$sheetExhibition.Range("F14").Value = 3419

$sheetAllTypes = $wb.Worksheets.Item("全部类型")
$sheetAllTypes.Range("F2").Value = 4350
$sheetAllTypes.Range("F13").Value = 157
$sheetAllTypes.Range("F16").Value = 1622
$sheetAllTypes.Range("F18").Value = 3419
